# Adds a new "Lens" line item to the "Buoy Sensors" sheet (row 33), which
# shifts every row from the old row 33 down through the old total row (48)
# down by one, updates the running-total formula/shared-formula range, and
# makes "Buoy Sensors" the active/selected sheet+cell (it was previously
# "3D Printed Buoy").

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("Buoy Sensors")

# --- Insert a new row at 33 (everything below, incl. formulas/hyperlinks,
#     shifts down by one row automatically) ------------------------------
$ws.Rows("33:33").Insert()

# --- Populate the new row with the Lens line item -----------------------
$ws.Range("A33").Value = "Lens"

$ws.Range("C33").Value = 35
$ws.Range("D33").Value = 1
$ws.Range("E33").Formula = "=C33*D33"

$ws.Range("F33").Value = 'Thorlabs - LA1540-AB N-BK7 Plano-Convex Lens, Ø1/2", f = 15 mm, AR Coating: 400 - 1100 nm'
$ws.Range("F33").Style = $ws.Range("F41").Style
$ws.Hyperlinks.Add($ws.Range("F33"), "https://www.thorlabs.com/thorproduct.cfm?partnumber=LA1540-AB") | Out-Null

# --- Make "Buoy Sensors" the active sheet/selection (previously "3D
#     Printed Buoy" was active) ------------------------------------------
$ws.Activate()
$ws.Range("E33").Select()
